# Highlight quantitative impact metrics (percentages, dollar amounts,
# large numbers) in bold + color (#2C3E50) across specific bullet points
# in the resume, matching the "hybrid bold + color highlighting" scheme
# described in the commit message.

$d = $word.ActiveDocument

# Color 2C3E50 expressed as a BGR-packed OLE_COLOR value (R + G*256 + B*65536)
$highlightColor = 5258796

function Highlight-InParagraph {
    param($ParaIndex, $Needles)

    $para = $d.Paragraphs.Item($ParaIndex)
    foreach ($needle in $Needles) {
        $rng = $para.Range.Duplicate
        $found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, `
                                    $true, 1, $false, "", 0)
        if ($found) {
            $rng.Font.Bold = $true
            $rng.Font.Color = $highlightColor
        }
    }
}

# "Partner - Siege Analytics" bullet: demographic classification accuracy
Highlight-InParagraph 9 @("23%", "64%")

# "Partner - Siege Analytics" bullet: turnout prediction / polling error margins
Highlight-InParagraph 11 @("87%", "71%", [char]0xB1 + "4.2%", [char]0xB1 + "2.1%")

# "Senior Analyst - Myers Research" bullet: RFP vendor bids
Highlight-InParagraph 31 @("1,200")

# "Programmer - Lake Research Partners" bullet: Polling Consortium Database value
Highlight-InParagraph 46 @("$400M", "$1B")

# "KEY ACHIEVEMENTS AND IMPACT" bullet: algorithm cost reduction
Highlight-InParagraph 63 @("73.5%", "$4.7M")

# "KEY ACHIEVEMENTS AND IMPACT" bullet: turnout prediction accuracy (short form)
Highlight-InParagraph 65 @("87%", "71%")
